$wb = $excel.ActiveWorkbook

# --- Update OTP (Open Targets Platform) source_version on the "compounds" sheet ---
$compounds = $wb.Worksheets.Item("compounds")
$compounds.Activate()
$compounds.Range("E2").Select()
$compounds.Range("E2").Value = "2023.12"

# --- Update Mitelman Database source_version on the "biomarkers" sheet ---
$biomarkers = $wb.Worksheets.Item("biomarkers")
$biomarkers.Activate()
$biomarkers.Range("E3").Select()
$biomarkers.Range("E3").Value = "v20231016"
